# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Black Amber" (Primera / Segunda)
# at the top of the existing block (old row 142), pushing the rest of
# the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 142-143 (existing rows 142.. shift down to 144..)
$ws.Range("A142:A143").EntireRow.Insert()

# New row 142: Black Amber / Primera
$ws.Cells.Item(142, 1).Value = 4
$ws.Cells.Item(142, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(142, 3).Value = "Los Lagos"
$ws.Cells.Item(142, 4).Value = 44617
$ws.Cells.Item(142, 5).Value = 10
$ws.Cells.Item(142, 6).Value = "Fruta"
$ws.Cells.Item(142, 7).Value = 100103
$ws.Cells.Item(142, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(142, 9).Value = 100103002
$ws.Cells.Item(142, 10).Value = "Ciruela"
$ws.Cells.Item(142, 11).Value = "Black Amber"
$ws.Cells.Item(142, 12).Value = "Primera"
$ws.Cells.Item(142, 13).Value = 400
$ws.Cells.Item(142, 14).Value = 15000
$ws.Cells.Item(142, 15).Value = 16000
$ws.Cells.Item(142, 16).Value = 15500
$ws.Cells.Item(142, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(142, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(142, 19).Value = 1033
$ws.Cells.Item(142, 20).Value = 15

# New row 143: Black Amber / Segunda
$ws.Cells.Item(143, 1).Value = 4
$ws.Cells.Item(143, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(143, 3).Value = "Los Lagos"
$ws.Cells.Item(143, 4).Value = 44617
$ws.Cells.Item(143, 5).Value = 10
$ws.Cells.Item(143, 6).Value = "Fruta"
$ws.Cells.Item(143, 7).Value = 100103
$ws.Cells.Item(143, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(143, 9).Value = 100103002
$ws.Cells.Item(143, 10).Value = "Ciruela"
$ws.Cells.Item(143, 11).Value = "Black Amber"
$ws.Cells.Item(143, 12).Value = "Segunda"
$ws.Cells.Item(143, 13).Value = 200
$ws.Cells.Item(143, 14).Value = 13000
$ws.Cells.Item(143, 15).Value = 13000
$ws.Cells.Item(143, 16).Value = 13000
$ws.Cells.Item(143, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(143, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(143, 19).Value = 867
$ws.Cells.Item(143, 20).Value = 15
